$d = $word.ActiveDocument

# 1. Replace all occurrences of "失败" (failed) with "失效" (failure) throughout
#    the document. This covers the heading "2. 失败管道的分析" -> "2. 失效管道的分析",
#    the sentence ending "...失败管道的基本情况见表1-1。" -> "...失效管道的基本情况见表1-1。",
#    and "在现场发现失败样品的穿孔部分..." -> "在现场发现失效样品的穿孔部分...".
$d.Content.Find.Execute("失败", $false, $false, $false, $false, $false, $true, 1, $false, "失效", 2) | Out-Null

# 2. Insert a new abstract-like paragraph right after the blank paragraph that
#    precedes the "2. 失效管道的分析" heading (i.e. directly before that heading).
$precedingBlank = $d.Paragraphs(4)
$precedingBlank.Range.InsertParagraphAfter()
$newPara = $precedingBlank.Next()

$fullText = "川南页岩气田已出现腐蚀破坏。 本文首先对腐蚀现象进行了分析。 其次，从腐蚀环境、材料、腐蚀产物等方面对腐蚀原因进行了深入分析。 最后对现场水环境进行模拟，结果表明：SRB是页岩气管道腐蚀穿孔的主要原因之一，细菌富集是局部腐蚀的主要原因，CO2促进细菌膜的形成并加速腐蚀。"
$newPara.Range.Text = $fullText

# Split the inserted text into two runs matching the source formatting split:
# first run holds the short lead-in sentence, second run holds the rest.
$firstSentence = "川南页岩气田已出现腐蚀破坏。"
$paraStart = $newPara.Range.Start
$splitPoint = $paraStart + $firstSentence.Length
$r1 = $d.Range($paraStart, $splitPoint)
$r1.Font.Size = $r1.Font.Size
